$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 249, shifting existing rows 249:266 down to 250:267
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row 249 with the new weekly record
$ws.Cells.Item(249, 1).Value = 4
$ws.Cells.Item(249, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(249, 3).Value = "Los Lagos"
$ws.Cells.Item(249, 4).Value = 44714
$ws.Cells.Item(249, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(249, 5).Value = 10
$ws.Cells.Item(249, 6).Value = 100112037
$ws.Cells.Item(249, 7).Value = "Cebollín"
$ws.Cells.Item(249, 8).Value = "Sin especificar"
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 10).Value = 70
$ws.Cells.Item(249, 11).Value = 11000
$ws.Cells.Item(249, 12).Value = 12000
$ws.Cells.Item(249, 13).Value = 11500
$ws.Cells.Item(249, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(249, 15).Value = "Región Metropolitana"
$ws.Cells.Item(249, 16).Value = 319
$ws.Cells.Item(249, 17).Value = 36
$ws.Cells.Item(249, 18).Value = "Hortaliza"
